$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.867.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -6.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.429.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -9.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.30"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.33"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -7.00%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0979"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -7.12%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.29"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.858.62"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -9.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.88"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -8.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.809.55"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.46%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.471.72"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -7.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -7.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.27"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.79"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.966"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.70"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -9.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.464"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.159"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.964"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.68"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.77"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.70%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0766"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -9.67%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.68"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.47"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.20"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.77%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.70"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.71%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "308.67"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -9.02%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.69"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.54"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.70"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.994"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.69"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0926"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.75%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.572"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.35%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0517"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.54%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.97"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -8.39%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.978.71"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.11%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0227"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.24%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.28"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -10.26%  "

$wb.Save()
